$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3 and 4 were (mistakenly) populated with different customer records;
# correct them to hold the same customer record as row 2 (columns C-H).
# Using Copy/PasteSpecial (values) instead of direct .Value assignment keeps
# the cells' existing (unstyled) formatting intact and avoids Excel's
# auto-detection turning numeric-looking / date-looking text into real
# numbers or dates.
$xlPasteValues = -4163

$ws.Range("C2:H2").Copy() | Out-Null
$ws.Range("C3:H3").PasteSpecial($xlPasteValues) | Out-Null
$ws.Range("C2:H2").Copy() | Out-Null
$ws.Range("C4:H4").PasteSpecial($xlPasteValues) | Out-Null

$excel.CutCopyMode = 0
